# Daily attendance processing - 2026-01-20 22:38:53
# Rotate the comma-separated "Recorded By" list in column G one position to
# the left (the first entry moves to the end) for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -eq "") { continue }

    $parts = $val -split ",\s*"

    if ($parts.Count -gt 1) {
        $trimmed = @()
        foreach ($p in $parts) { $trimmed += $p.Trim() }

        $rotated = @()
        for ($i = 1; $i -lt $trimmed.Count; $i++) { $rotated += $trimmed[$i] }
        $rotated += $trimmed[0]

        $cell.Value = [string]::Join(", ", $rotated)
    }
}
